$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns with numeric-looking text values keep their original text type
# (Antal=I, Starttid=Z, Sluttid=AB, Publik kommentar=AC already default to General,
# but purely-numeric-looking strings like "10" or time-looking "07:30" could be
# auto-converted by Excel; set NumberFormat to Text for the affected cells to be safe.)

$ws.Range("A32").Value = 112204281
$ws.Range("Q32").Value = 595169
$ws.Range("R32").Value = 6396054
$ws.Range("A33").Value = 112203709
$ws.Range("B33").Value = 89082
$ws.Range("D33").Value = 'LC'
$ws.Range("E33").Value = 5741
$ws.Range("F33").Value = 'Tjockfotad fingersvamp'
$ws.Range("G33").Value = 'Ramaria flavescens'
$ws.Range("H33").Value = '(Schaeff.) R. H. Petersen'
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = '4'
$ws.Range("P33").Value = 'A 30779, Storön, Sm'
$ws.Range("Q33").Value = 594781
$ws.Range("R33").Value = 6396169
$ws.Range("S33").Value = 10
$ws.Range("AC33").Value = ""
$ws.Range("A35").Value = 112204311
$ws.Range("B35").Value = 96735
$ws.Range("D35").Value = 'VU'
$ws.Range("E35").Value = 220787
$ws.Range("F35").Value = 'Knärot'
$ws.Range("G35").Value = 'Goodyera repens'
$ws.Range("H35").Value = '(L.) R. Br.'
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = '300'
$ws.Range("J35").Value = 'plantor/tuvor'
$ws.Range("P35").Value = 'Storön, Samsvik, Sm'
$ws.Range("Q35").Value = 595112
$ws.Range("R35").Value = 6396025
$ws.Range("S35").Value = 5
$ws.Range("AW35").Value = 'Larsgunnar Nilsson'
$ws.Range("AX35").Value = 'Larsgunnar Nilsson, Gunilla Nilsson, Ingvor Kasselstrand, Magnus Kasselstrand'
$ws.Range("A36").Value = 112203732
$ws.Range("B36").Value = 90806
$ws.Range("D36").Value = 'NT'
$ws.Range("E36").Value = 4361
$ws.Range("F36").Value = 'Orange taggsvamp'
$ws.Range("G36").Value = 'Hydnellum aurantiacum'
$ws.Range("H36").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = '30'
$ws.Range("J36").Value = 'fruktkroppar'
$ws.Range("P36").Value = 'A 30779, Storön, Sm'
$ws.Range("Q36").Value = 594803
$ws.Range("R36").Value = 6396141
$ws.Range("S36").Value = 10
$ws.Range("AW36").Value = 'Magnus Kasselstrand'
$ws.Range("AX36").Value = 'Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson'
$ws.Range("A37").Value = 112204167
$ws.Range("B37").Value = 90837
$ws.Range("D37").Value = 'NT'
$ws.Range("E37").Value = 5966
$ws.Range("F37").Value = 'Motaggsvamp'
$ws.Range("G37").Value = 'Sarcodon squamosus'
$ws.Range("H37").Value = '(Schaeff.) Quél.'
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = '10'
$ws.Range("J37").Value = 'fruktkroppar'
$ws.Range("P37").Value = 'A 30799, Storön, Sm'
$ws.Range("Q37").Value = 594925
$ws.Range("R37").Value = 6396228
$ws.Range("S37").Value = 100
$ws.Range("AC37").Value = '2+7+1'
$ws.Range("AW37").Value = 'Magnus Kasselstrand'
$ws.Range("AX37").Value = 'Magnus Kasselstrand, Ingvor Kasselstrand, Gunilla Nilsson, Larsgunnar Nilsson'
$ws.Range("A38").Value = 112204297
$ws.Range("B38").Value = 96735
$ws.Range("D38").Value = 'VU'
$ws.Range("E38").Value = 220787
$ws.Range("F38").Value = 'Knärot'
$ws.Range("G38").Value = 'Goodyera repens'
$ws.Range("H38").Value = '(L.) R. Br.'
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = '10'
$ws.Range("J38").Value = 'plantor/tuvor'
$ws.Range("P38").Value = 'Storön, Samsvik, Sm'
$ws.Range("Q38").Value = 595097
$ws.Range("R38").Value = 6396058
$ws.Range("S38").Value = 5
$ws.Range("AW38").Value = 'Larsgunnar Nilsson'
$ws.Range("AX38").Value = 'Larsgunnar Nilsson, Gunilla Nilsson, Ingvor Kasselstrand, Magnus Kasselstrand'
$ws.Range("A39").Value = 112236300
$ws.Range("B39").Value = 103781
$ws.Range("E39").Value = 221144
$ws.Range("F39").Value = 'Grönpyrola'
$ws.Range("G39").Value = 'Pyrola chlorantha'
$ws.Range("H39").Value = 'Sw.'
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = '30'
$ws.Range("J39").Value = 'plantor/tuvor'
$ws.Range("Q39").Value = 594877
$ws.Range("R39").Value = 6396173
$ws.Range("A40").Value = 112236205
$ws.Range("B40").Value = 89573
$ws.Range("D40").Value = 'NT'
$ws.Range("E40").Value = 5442
$ws.Range("F40").Value = 'Tallticka'
$ws.Range("G40").Value = 'Porodaedalea pini'
$ws.Range("H40").Value = '(Brot.) Murrill'
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = '1'
$ws.Range("J40").Value = 'fruktkroppar'
$ws.Range("Q40").Value = 594931
$ws.Range("R40").Value = 6396214
$ws.Range("A43").Value = 112236185
$ws.Range("B43").Value = 96735
$ws.Range("D43").Value = 'VU'
$ws.Range("E43").Value = 220787
$ws.Range("F43").Value = 'Knärot'
$ws.Range("G43").Value = 'Goodyera repens'
$ws.Range("H43").Value = '(L.) R. Br.'
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = '300'
$ws.Range("J43").Value = 'plantor/tuvor'
$ws.Range("P43").Value = 'A 30779, Storön, Sm'
$ws.Range("Q43").Value = 594982
$ws.Range("R43").Value = 6396167
$ws.Range("S43").Value = 25
$ws.Range("Z43").Value = '07:30'
$ws.Range("AB43").Value = '10:45'
$ws.Range("AC43").Value = ""
$ws.Range("A44").Value = 112236402
$ws.Range("B44").Value = 103781
$ws.Range("D44").Value = 'LC'
$ws.Range("E44").Value = 221144
$ws.Range("F44").Value = 'Grönpyrola'
$ws.Range("G44").Value = 'Pyrola chlorantha'
$ws.Range("H44").Value = 'Sw.'
$ws.Range("I44").NumberFormat = "@"
$ws.Range("I44").Value = '50'
$ws.Range("J44").Value = 'plantor/tuvor'
$ws.Range("Q44").Value = 594889
$ws.Range("R44").Value = 6396160
$ws.Range("A45").Value = 112236282
$ws.Range("B45").Value = 96735
$ws.Range("D45").Value = 'VU'
$ws.Range("E45").Value = 220787
$ws.Range("F45").Value = 'Knärot'
$ws.Range("G45").Value = 'Goodyera repens'
$ws.Range("H45").Value = '(L.) R. Br.'
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = '10'
$ws.Range("Q45").Value = 594909
$ws.Range("R45").Value = 6396198
$ws.Range("A46").Value = 112236343
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = '10'
$ws.Range("Q46").Value = 594869
$ws.Range("R46").Value = 6396235
$ws.Range("A47").Value = 112237682
$ws.Range("B47").Value = 90335
$ws.Range("D47").Value = 'NT'
$ws.Range("E47").Value = 2014
$ws.Range("F47").Value = 'Koralltaggsvamp'
$ws.Range("G47").Value = 'Hericium coralloides'
$ws.Range("H47").Value = '(Scop.:Fr.) Pers.'
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value = '1'
$ws.Range("J47").Value = 'fruktkroppar'
$ws.Range("P47").Value = 'Storön, Sm'
$ws.Range("Q47").Value = 595194
$ws.Range("R47").Value = 6396077
$ws.Range("Z47").Value = ""
$ws.Range("AB47").Value = ""
$ws.Range("AC47").Value = 'I granplantering'
$ws.Range("A48").Value = 112236468
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = '20'
$ws.Range("Q48").Value = 594918
$ws.Range("R48").Value = 6396098
$ws.Range("A49").Value = 112236222
$ws.Range("B49").Value = 89950
$ws.Range("D49").Value = 'LC'
$ws.Range("E49").Value = 5420
$ws.Range("F49").Value = 'Grovticka'
$ws.Range("G49").Value = 'Phaeolus schweinitzii'
$ws.Range("H49").Value = '(Fr.) Pat.'
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value = '1'
$ws.Range("J49").Value = 'mycel'
$ws.Range("Q49").Value = 594933
$ws.Range("R49").Value = 6396201
$ws.Range("S49").Value = 10
